{"js": "// Insert a brand-new right-aligned, bold \"Times New Roman\" 14pt paragraph\n// (\"\u041a\u0443\u0434\u0438\u043d\u043e\u0432 \u0412\u043b\u0430\u0434\u0438\u0441\u043b\u0430\u0432 \u041f\u042023/2\") as the very first paragraph of the document\n// body, with 1.5 (360 twips \"auto\") line spacing \u2014 matching the author's\n// heading/name-block style used throughout the rest of the document.\nconst body = context.document.body;\n\nconst namePara = body.insertParagraph(\n  \"\u041a\u0443\u0434\u0438\u043d\u043e\u0432 \u0412\u043b\u0430\u0434\u0438\u0441\u043b\u0430\u0432 \u041f\u042023/2\",\n  Word.InsertLocation.start\n);\n\nnamePara.alignment = Word.Alignment.right;\nnamePara.lineSpacing = 18; // 360 twentieths-of-a-point == 18pt == 1.5 lines\n\nnamePara.font.name = \"Times New Roman\";\nnamePara.font.bold = true;\nnamePara.font.size = 14; // half-points 28 == 14pt\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Insert a brand-new empty paragraph immediately before the current first\n# paragraph of the document body (at the very start of the main story).\n$r = $d.Range(0, 0)\n$r.InsertParagraphBefore() | Out-Null\n\n# The freshly inserted paragraph is now paragraph #1 - fill in its text and\n# apply the same right-aligned, bold \"Times New Roman\" 14pt, 1.5-line-spacing\n# formatting used for the other heading lines in this document.\n$p = $d.Paragraphs(1)\n$p.Range.Text = \"\u041a\u0443\u0434\u0438\u043d\u043e\u0432 \u0412\u043b\u0430\u0434\u0438\u0441\u043b\u0430\u0432 \u041f\u042023/2\"\n\n$p.Alignment = \"wdAlignParagraphRight\"\n$p.Range.ParagraphFormat.LineSpacingRule = \"wdLineSpace1pt5\"\n\n$p.Range.Font.Name = \"Times New Roman\"\n$p.Range.Font.Bold = $true\n$p.Range.Font.Size = 14\n"}
